$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = "diya"
$ws.Range("C3").Value = "cs22B"
$ws.Range("D3").Value = "86h"
$ws.Range("E3").Value = "diya"
$ws.Range("F3").Value = (Get-Date -Year 1899 -Month 12 -Day 30 -Hour 0 -Minute 0 -Second 0).AddDays(43594)
$ws.Range("G3").Value = (Get-Date -Year 1899 -Month 12 -Day 30 -Hour 0 -Minute 0 -Second 0).AddDays(43607)
$ws.Range("H3").Value = (Get-Date -Year 1899 -Month 12 -Day 30 -Hour 0 -Minute 0 -Second 0).AddDays(43607)
